$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Address for the last row (Aldermore Bank Plc) which previously
# contained a stray phone number instead of the actual address.
$ws.Range("F7").Value = "50 St Mary Axe, KT17 1BS, London, United Kingdom"
